# Auto-generated Excel COM-interop script
# Applies the quarterly/annual 2022-2025 forecast extension described in the commit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# INDI: append quarterly rows 106-121 (DATE / INDI_NB / IND_base)
# ---------------------------------------------------------------
$wsINDI = $wb.Worksheets.Item("INDI")
$indiDates = $wsINDI.Range("A106:A121")
$indiDates.NumberFormat = "@"
$indiData = @(
    @("2022-01-01", 649168.848173915, 64.9168848173915),
    @("2022-04-01", 1335507.88342838, 133.550788342838),
    @("2022-07-01", 1328672.5794282, 132.86725794282),
    @("2022-10-01", 1445249.48223089, 144.524948223089),
    @("2023-01-01", 672197.43653263, 67.219743653263),
    @("2023-04-01", 1358536.4717871, 135.85364717871),
    @("2023-07-01", 1351701.16778691, 135.170116778691),
    @("2023-10-01", 1468278.0705896, 146.82780705896),
    @("2024-01-01", 695226.024891344, 69.5226024891344),
    @("2024-04-01", 1381565.06014581, 138.156506014581),
    @("2024-07-01", 1374729.75614563, 137.472975614563),
    @("2024-10-01", 1491306.65894831, 149.130665894831),
    @("2025-01-01", 718254.613250059, 71.8254613250059),
    @("2025-04-01", 1404593.64850453, 140.459364850453),
    @("2025-07-01", 1397758.34450434, 139.775834450434),
    @("2025-10-01", 1514335.24730703, 151.433524730703),
)
for ($i = 0; $i -lt $indiData.Length; $i++) {
    $r = 106 + $i
    $row = $indiData[$i]
    $wsINDI.Cells.Item($r, 1).Value = $row[0]
    $wsINDI.Cells.Item($r, 2).Value = $row[1]
    $wsINDI.Cells.Item($r, 3).Value = $row[2]
}
$indiDates.Style = "Normal"

# ---------------------------------------------------------------
# ETALONNAGE: append annual rows 28-31
# ---------------------------------------------------------------
$wsETAL = $wb.Worksheets.Item("ETALONNAGE")
$etalTextA = $wsETAL.Range("A28:A31")
$etalTextC = $wsETAL.Range("C28:C31")
$etalTextF = $wsETAL.Range("F28:F31")
$etalTextA.NumberFormat = "@"
$etalTextC.NumberFormat = "@"
$etalTextF.NumberFormat = "@"
$etalData = @(
    @("2022", 2296616.87359465, "118.964969831534", 1.42901178624326, 1.83570249323546, "Acceptable"),
    @("2023", 2331363.1405257, "121.267828667406", 1.51293266763566, 1.93574531993117, "Acceptable"),
    @("2024", 2366109.40745675, "123.570687503277", 1.49038415882357, 1.89898579134897, "Acceptable"),
    @("2025", 2400855.67438779, "125.873546339149", 1.46849789876768, 1.86359636124089, "Acceptable"),
)
for ($i = 0; $i -lt $etalData.Length; $i++) {
    $r = 28 + $i
    $row = $etalData[$i]
    $wsETAL.Cells.Item($r, 1).Value = $row[0]
    $wsETAL.Cells.Item($r, 2).Value = $row[1]
    $wsETAL.Cells.Item($r, 3).Value = $row[2]
    $wsETAL.Cells.Item($r, 4).Value = $row[3]
    $wsETAL.Cells.Item($r, 5).Value = $row[4]
    $wsETAL.Cells.Item($r, 6).Value = $row[5]
}
$etalTextA.Style = "Normal"
$etalTextC.Style = "Normal"
$etalTextF.Style = "Normal"

# ---------------------------------------------------------------
# PREVISION: append quarterly rows 102-117
# ---------------------------------------------------------------
$wsPREV = $wb.Worksheets.Item("PREVISION")
$prevTextA = $wsPREV.Range("A102:A117")
$prevTextK = $wsPREV.Range("K102:K117")
$prevTextA.NumberFormat = "@"
$prevTextK.NumberFormat = "@"
$prevData = @(
    @("2022-01-01", 370280.541931793, 239321.877687495, 233178.042582547, 137102.499349247, 16.2292212043479, 66460.5912471754, 14746.3562591272, 3.49131025750811, 2.28186969845432, "Acceptable"),
    @("2022-04-01", 629173.086197585, 400613.685764224, 396211.337273556, 232961.74892403, 33.3876970857095, 66460.5912471754, 11998.8415114648, 1.48320460545033, 1.18406675006326, "Acceptable"),
    @("2022-07-01", 626594.755150478, 400262.752678334, 394587.675971891, 232007.079178588, 33.2168144857049, 66460.5912471754, 12050.0041583033, 1.63169114736716, 1.3008704324762, "Acceptable"),
    @("2022-10-01", 670568.490314797, 428567.21831726, 422279.407860271, 248289.082454527, 36.1312370557721, 66460.5912471754, 11861.4045142082, 1.61919913129787, 1.31240428591117, "Acceptable"),
    @("2023-01-01", 378967.108664555, 247811.571889555, 238648.26420678, 140318.844457775, 16.8049359133157, 66460.5912471754, 14746.3562591272, 3.54739578516328, 2.34594199507296, "Acceptable"),
    @("2023-04-01", 637859.652930347, 407521.595312949, 401681.558897789, 236178.094032558, 33.9634117946774, 66460.5912471754, 11998.8415114648, 1.72433189234329, 1.38063228121519, "Acceptable"),
    @("2023-07-01", 635281.32188324, 407200.117315393, 400057.897596124, 235223.424287116, 33.7925291946728, 66460.5912471754, 12050.0041583033, 1.7332026501687, 1.38631334867707, "Acceptable"),
    @("2023-10-01", 679255.057047559, 435396.003365108, 427749.629484503, 251505.427563055, 36.70695176474, 66460.5912471754, 11861.4045142082, 1.59339883126388, 1.29540335673743, "Acceptable"),
    @("2024-01-01", 387653.675397317, 256301.266091615, 244118.485831013, 143535.189566304, 17.3806506222836, 66460.5912471754, 14746.3562591272, 3.42586673306908, 2.29216903898923, "Acceptable"),
    @("2024-04-01", 646546.219663109, 414429.504861673, 407151.780522022, 239394.439141087, 34.5391265036453, 66460.5912471754, 11998.8415114648, 1.69510269594906, 1.36183041094624, "Acceptable"),
    @("2024-07-01", 643967.888616002, 414137.481952452, 405528.119220357, 238439.769395645, 34.3682439036406, 66460.5912471754, 12050.0041583033, 1.70367451826785, 1.36735748927912, "Acceptable"),
    @("2024-10-01", 687941.62378032, 442224.788412955, 433219.851108736, 254721.772671584, 37.2826664737079, 66460.5912471754, 11861.4045142082, 1.56840783908645, 1.27883725599609, "Acceptable"),
    @("2025-01-01", 396340.242130078, 264790.960293675, 249588.707455246, 146751.534674833, 17.9563653312515, 66460.5912471754, 14746.3562591272, 3.31238871017718, 2.24080597813463, "Acceptable"),
    @("2025-04-01", 655232.78639587, 421337.414410397, 412622.002146255, 242610.784249615, 35.1148412126131, 66460.5912471754, 11998.8415114648, 1.66684791205445, 1.34353375962633, "Acceptable"),
    @("2025-07-01", 652654.455348763, 421074.846589511, 410998.34084459, 241656.114504173, 34.9439586126085, 66460.5912471754, 12050.0041583033, 1.67513565890074, 1.34891302599427, "Acceptable"),
    @("2025-10-01", 696628.190513082, 449053.573460803, 438690.07273297, 257938.117780113, 37.8583811826757, 66460.5912471754, 11861.4045142082, 1.54418866304495, 1.26268951208794, "Acceptable"),
)
for ($i = 0; $i -lt $prevData.Length; $i++) {
    $r = 102 + $i
    $row = $prevData[$i]
    $wsPREV.Cells.Item($r, 1).Value = $row[0]
    $wsPREV.Cells.Item($r, 2).Value = $row[1]
    $wsPREV.Cells.Item($r, 3).Value = $row[2]
    $wsPREV.Cells.Item($r, 4).Value = $row[3]
    $wsPREV.Cells.Item($r, 5).Value = $row[4]
    $wsPREV.Cells.Item($r, 6).Value = $row[5]
    $wsPREV.Cells.Item($r, 7).Value = $row[6]
    $wsPREV.Cells.Item($r, 8).Value = $row[7]
    $wsPREV.Cells.Item($r, 9).Value = $row[8]
    $wsPREV.Cells.Item($r, 10).Value = $row[9]
    $wsPREV.Cells.Item($r, 11).Value = $row[10]
}
$prevTextA.Style = "Normal"
$prevTextK.Style = "Normal"

# ---------------------------------------------------------------
# VATRIM: append quarterly rows 102-117
# ---------------------------------------------------------------
$wsVAT = $wb.Worksheets.Item("VATRIM")
$vatText = $wsVAT.Range("A102:A117")
$vatText.NumberFormat = "@"
$vatData = @(
    @("2022-01-01", 239321.877687495),
    @("2022-04-01", 400613.685764224),
    @("2022-07-01", 400262.752678334),
    @("2022-10-01", 428567.21831726),
    @("2023-01-01", 247811.571889555),
    @("2023-04-01", 407521.595312949),
    @("2023-07-01", 407200.117315393),
    @("2023-10-01", 435396.003365108),
    @("2024-01-01", 256301.266091615),
    @("2024-04-01", 414429.504861673),
    @("2024-07-01", 414137.481952452),
    @("2024-10-01", 442224.788412955),
    @("2025-01-01", 264790.960293675),
    @("2025-04-01", 421337.414410397),
    @("2025-07-01", 421074.846589511),
    @("2025-10-01", 449053.573460803),
)
for ($i = 0; $i -lt $vatData.Length; $i++) {
    $r = 102 + $i
    $row = $vatData[$i]
    $wsVAT.Cells.Item($r, 1).Value = $row[0]
    $wsVAT.Cells.Item($r, 2).Value = $row[1]
}
$vatText.Style = "Normal"

